# Insert a new weekly data row before the existing row 703, shifting all
# subsequent rows down by one (754 -> 755). Then populate the new row
# with the latest weekly price entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 703; existing rows 703:754 shift to 704:755
$ws.Rows.Item(703).Insert()

# Populate the newly inserted row 703 with the new weekly record
$ws.Range("A703").Value = 3
$ws.Range("B703").Value = "Femacal de La Calera"
$ws.Range("C703").Value = "Coquimbo"
$ws.Range("D703").Value = 45265
$ws.Range("E703").Value = 5
$ws.Range("F703").Value = "Fruta"
$ws.Range("G703").Value = 100108
$ws.Range("H703").Value = "Tropicales y subtropicales"
$ws.Range("I703").Value = 100108002
$ws.Range("J703").Value = "Mango"
$ws.Range("K703").Value = "Sin especificar"
$ws.Range("L703").Value = "Primera"
$ws.Range("M703").Value = 228
$ws.Range("N703").Value = 12000
$ws.Range("O703").Value = 12000
$ws.Range("P703").Value = 12000
$ws.Range("Q703").Value = "$/bandeja 4 kilos"
$ws.Range("R703").Value = "Perú"
$ws.Range("S703").Value = 3000
$ws.Range("T703").Value = 4
